# Rename the sheet/tab: "Producto" -> "Inventario"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Inventario"

# Row 2: NOMBRE / value -> "asdfsdfs"
$ws.Range("B2").Value = "asdfsdfs"

# Row 3: label -> STOCK ACTUAL, value -> 0 (numeric)
$ws.Range("A3").Value = "STOCK ACTUAL"
$ws.Range("B3").Value = 0

# Row 4: label -> MÁXIMO STOCK, value -> 0 (numeric)
$ws.Range("A4").Value = "MÁXIMO STOCK"
$ws.Range("B4").Value = 0

# Row 5: label -> FECHA DE ACTUALIZACIÓN, value -> "22/11/2024" (text, not date)
$ws.Range("A5").Value = "FECHA DE ACTUALIZACIÓN"
$ws.Range("B5").Value = "22/11/2024"

# Remove old rows 6-9 entirely (no longer part of the report)
$ws.Range("A6:B9").EntireRow.Delete()
